$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, $CellRef, $Text)
    $Sheet.Range($CellRef).Value = "'" + $Text
    $Sheet.Range($CellRef).Style = "Normal"
}

# Row 2
Set-CellText $ws "D2" '67.499.23'
Set-CellText $ws "E2" '  -0.94%  '

# Row 3
Set-CellText $ws "D3" '3.235.95'
Set-CellText $ws "E3" '  -1.18%  '

# Row 4
Set-CellText $ws "E4" '  +0.00%  '

# Row 5
Set-CellText $ws "D5" '578.99'
Set-CellText $ws "E5" '  -1.55%  '

# Row 6
Set-CellText $ws "D6" '183.76'
Set-CellText $ws "E6" '  -1.44%  '

# Row 7
Set-CellText $ws "B7" 'USDC'
Set-CellText $ws "C7" 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-CellText $ws "D7" '1.00'
Set-CellText $ws "E7" '  +0.00%  '

# Row 8
Set-CellText $ws "B8" 'XRP'
Set-CellText $ws "C8" 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-CellText $ws "D8" '0.609'
Set-CellText $ws "E8" '  +1.33%  '

# Row 9
Set-CellText $ws "D9" '3.234.18'
Set-CellText $ws "E9" '  -1.22%  '

# Row 10
Set-CellText $ws "D10" '0.131'
Set-CellText $ws "E10" '  -3.29%  '

# Row 11
Set-CellText $ws "E11" '  -2.29%  '

# Row 13
Set-CellText $ws "D13" '3.792.41'
Set-CellText $ws "E13" '  -1.34%  '

# Row 14
Set-CellText $ws "E14" '  +0.19%  '

# Row 15
Set-CellText $ws "D15" '27.63'
Set-CellText $ws "E15" '  -3.59%  '

# Row 16
Set-CellText $ws "D16" '67.542.57'
Set-CellText $ws "E16" '  -0.90%  '

# Row 17
Set-CellText $ws "E17" '  -1.97%  '

# Row 18
Set-CellText $ws "D18" '3.221.62'
Set-CellText $ws "E18" '  -1.79%  '

# Row 19
Set-CellText $ws "D19" '5.75'
Set-CellText $ws "E19" '  -1.76%  '

# Row 20
Set-CellText $ws "E20" '  -1.17%  '

# Row 21
Set-CellText $ws "E21" '  +3.21%  '

# Row 22
Set-CellText $ws "D22" '7.56'
Set-CellText $ws "E22" '  -2.22%  '

# Row 23
Set-CellText $ws "E23" '  -0.09%  '

# Row 25
Set-CellText $ws "E25" '  -0.07%  '

# Row 26
Set-CellText $ws "E26" '  -2.56%  '

# Row 27
Set-CellText $ws "D27" '0.187'
Set-CellText $ws "E27" '  -1.33%  '

# Row 28
Set-CellText $ws "D28" '9.63'
Set-CellText $ws "E28" '  -1.63%  '

# Row 29
Set-CellText $ws "D29" '1.00'
Set-CellText $ws "E29" '  +0.16%  '

# Row 30
Set-CellText $ws "D30" '1.95'
Set-CellText $ws "E30" '  -2.18%  '

# Row 31
Set-CellText $ws "D31" '5.55'
Set-CellText $ws "E31" '  -4.43%  '

# Row 32
Set-CellText $ws "D32" '22.59'
Set-CellText $ws "E32" '  -1.56%  '

# Row 33
Set-CellText $ws "D33" '7.00'
Set-CellText $ws "E33" '  -2.21%  '

# Row 34
Set-CellText $ws "E34" '  -2.35%  '

# Row 36
Set-CellText $ws "D36" '160.61'
Set-CellText $ws "E36" '  -1.46%  '

# Row 37
Set-CellText $ws "E37" '  -4.13%  '

# Row 38
Set-CellText $ws "E38" '  +0.82%  '

# Row 39
Set-CellText $ws "E39" '  -0.93%  '

# Row 40
Set-CellText $ws "D40" '0.802'
Set-CellText $ws "E40" '  -4.46%  '

# Row 41
Set-CellText $ws "E41" '  -1.54%  '

# Row 42
Set-CellText $ws "E42" '  -4.84%  '

# Row 43
Set-CellText $ws "E43" '  -6.31%  '

# Row 44
Set-CellText $ws "E44" '  -0.70%  '

# Row 45
Set-CellText $ws "B45" 'Maker'
Set-CellText $ws "C45" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws "D45" '2.612.37'
Set-CellText $ws "E45" '  -1.26%  '

# Row 46
Set-CellText $ws "B46" 'OKB'
Set-CellText $ws "C46" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText $ws "D46" '40.51'
Set-CellText $ws "E46" '  -2.06%  '

# Row 47
Set-CellText $ws "D47" '24.70'
Set-CellText $ws "E47" '  -3.08%  '

# Row 48
Set-CellText $ws "D48" '333.91'
Set-CellText $ws "E48" '  -3.20%  '

# Row 49
Set-CellText $ws "E49" '  -2.16%  '

# Row 50
Set-CellText $ws "D50" '6.32'
Set-CellText $ws "E50" '  +0.82%  '

# Row 51
Set-CellText $ws "E51" '  -0.85%  '
